$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 131046822
$ws.Range("B8").Value = 79243
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("P8").Value = "Blomkällan, Dlr"
$ws.Range("Q8").Value = 401699
$ws.Range("R8").Value = 6818070
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = "Dalarna"
$ws.Range("U8").Value = "Älvdalen"
$ws.Range("V8").Value = "Dalarna"
$ws.Range("W8").Value = "Särna"
$ws.Range("Y8").Value = "'2026-02-05"
$ws.Range("Z8").Value = "'14:49"
$ws.Range("AA8").Value = "'2026-02-05"
$ws.Range("AB8").Value = "'14:49"
$ws.Range("AD8").Value = $False
$ws.Range("AE8").Value = $False
$ws.Range("AG8").Value = $False
$ws.Range("AW8").Value = "Philipp Weiss"
$ws.Range("AX8").Value = "Philipp Weiss"
$ws.Range("M8").ClearContents()
$ws.Range("AC8").ClearContents()

# Row 9
$ws.Range("A9").Value = 131047013
$ws.Range("B9").Value = 57884
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("M9").Value = "färska spår"
$ws.Range("P9").Value = "Blomkällan, Dlr"
$ws.Range("Q9").Value = 401631
$ws.Range("R9").Value = 6817903
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Dalarna"
$ws.Range("U9").Value = "Älvdalen"
$ws.Range("V9").Value = "Dalarna"
$ws.Range("W9").Value = "Särna"
$ws.Range("Y9").Value = "'2026-02-05"
$ws.Range("Z9").Value = "'14:57"
$ws.Range("AA9").Value = "'2026-02-05"
$ws.Range("AB9").Value = "'14:57"
$ws.Range("AC9").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AD9").Value = $False
$ws.Range("AE9").Value = $True
$ws.Range("AG9").Value = $False
$ws.Range("AW9").Value = "Philipp Weiss"
$ws.Range("AX9").Value = "Philipp Weiss"

# Row 10
$ws.Range("A10").Value = 131046773
$ws.Range("B10").Value = 57884
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("M10").Value = "äldre spår"
$ws.Range("P10").Value = "Blomkällan, Dlr"
$ws.Range("Q10").Value = 401346
$ws.Range("R10").Value = 6818162
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = "Dalarna"
$ws.Range("U10").Value = "Älvdalen"
$ws.Range("V10").Value = "Dalarna"
$ws.Range("W10").Value = "Särna"
$ws.Range("Y10").Value = "'2026-02-05"
$ws.Range("Z10").Value = "'15:23"
$ws.Range("AA10").Value = "'2026-02-05"
$ws.Range("AB10").Value = "'15:23"
$ws.Range("AC10").Value = "Äldre ringhack (gran)"
$ws.Range("AD10").Value = $False
$ws.Range("AE10").Value = $False
$ws.Range("AG10").Value = $False
$ws.Range("AW10").Value = "Philipp Weiss"
$ws.Range("AX10").Value = "Philipp Weiss"

# Row 11
$ws.Range("A11").Value = 131046823
$ws.Range("B11").Value = 79243
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("P11").Value = "Blomkällan, Dlr"
$ws.Range("Q11").Value = 401661
$ws.Range("R11").Value = 6818064
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = "Dalarna"
$ws.Range("U11").Value = "Älvdalen"
$ws.Range("V11").Value = "Dalarna"
$ws.Range("W11").Value = "Särna"
$ws.Range("Y11").Value = "'2026-02-05"
$ws.Range("Z11").Value = "'14:50"
$ws.Range("AA11").Value = "'2026-02-05"
$ws.Range("AB11").Value = "'14:50"
$ws.Range("AD11").Value = $False
$ws.Range("AE11").Value = $False
$ws.Range("AG11").Value = $False
$ws.Range("AW11").Value = "Philipp Weiss"
$ws.Range("AX11").Value = "Philipp Weiss"
$ws.Range("M11").ClearContents()
$ws.Range("AC11").ClearContents()

# Row 16
$ws.Range("A16").Value = 131046708
$ws.Range("B16").Value = 83223
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 6440
$ws.Range("F16").Value = "Vitgrynig nållav"
$ws.Range("G16").Value = "Chaenotheca subroscida"
$ws.Range("H16").Value = "(Eitner) Zahlbr."
$ws.Range("P16").Value = "Blomkällan, Dlr"
$ws.Range("Q16").Value = 401645
$ws.Range("R16").Value = 6818016
$ws.Range("S16").Value = 10
$ws.Range("T16").Value = "Dalarna"
$ws.Range("U16").Value = "Älvdalen"
$ws.Range("V16").Value = "Dalarna"
$ws.Range("W16").Value = "Särna"
$ws.Range("Y16").Value = "'2026-02-05"
$ws.Range("Z16").Value = "'14:52"
$ws.Range("AA16").Value = "'2026-02-05"
$ws.Range("AB16").Value = "'14:52"
$ws.Range("AD16").Value = $False
$ws.Range("AE16").Value = $False
$ws.Range("AG16").Value = $False
$ws.Range("AW16").Value = "Philipp Weiss"
$ws.Range("AX16").Value = "Philipp Weiss"

# Row 17
$ws.Range("A17").Value = 131046724
$ws.Range("B17").Value = 79275
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 185
$ws.Range("F17").Value = "Violettgrå tagellav"
$ws.Range("G17").Value = "Bryoria nadvornikiana"
$ws.Range("H17").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("P17").Value = "Blomkällan, Dlr"
$ws.Range("Q17").Value = 401635
$ws.Range("R17").Value = 6817874
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = "Dalarna"
$ws.Range("U17").Value = "Älvdalen"
$ws.Range("V17").Value = "Dalarna"
$ws.Range("W17").Value = "Särna"
$ws.Range("Y17").Value = "'2026-02-05"
$ws.Range("Z17").Value = "'14:58"
$ws.Range("AA17").Value = "'2026-02-05"
$ws.Range("AB17").Value = "'14:58"
$ws.Range("AD17").Value = $False
$ws.Range("AE17").Value = $False
$ws.Range("AG17").Value = $False
$ws.Range("AW17").Value = "Philipp Weiss"
$ws.Range("AX17").Value = "Philipp Weiss"

# Row 18
$ws.Range("A18").Value = 131046768
$ws.Range("B18").Value = 57884
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("M18").Value = "färska spår"
$ws.Range("P18").Value = "Blomkällan, Dlr"
$ws.Range("Q18").Value = 401346
$ws.Range("R18").Value = 6818439
$ws.Range("S18").Value = 10
$ws.Range("T18").Value = "Dalarna"
$ws.Range("U18").Value = "Älvdalen"
$ws.Range("V18").Value = "Dalarna"
$ws.Range("W18").Value = "Särna"
$ws.Range("Y18").Value = "'2026-02-05"
$ws.Range("Z18").Value = "'15:32"
$ws.Range("AA18").Value = "'2026-02-05"
$ws.Range("AB18").Value = "'15:32"
$ws.Range("AC18").Value = "Både färska och gamla ringhack på grov gammal tall"
$ws.Range("AD18").Value = $False
$ws.Range("AE18").Value = $False
$ws.Range("AG18").Value = $False
$ws.Range("AW18").Value = "Philipp Weiss"
$ws.Range("AX18").Value = "Philipp Weiss"

# Row 19
$ws.Range("A19").Value = 131046827
$ws.Range("B19").Value = 79243
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = "Garnlav"
$ws.Range("G19").Value = "Alectoria sarmentosa"
$ws.Range("H19").Value = "(Ach.) Ach."
$ws.Range("P19").Value = "Blomkällan, Dlr"
$ws.Range("Q19").Value = 401630
$ws.Range("R19").Value = 6817911
$ws.Range("S19").Value = 10
$ws.Range("T19").Value = "Dalarna"
$ws.Range("U19").Value = "Älvdalen"
$ws.Range("V19").Value = "Dalarna"
$ws.Range("W19").Value = "Särna"
$ws.Range("Y19").Value = "'2026-02-05"
$ws.Range("Z19").Value = "'14:56"
$ws.Range("AA19").Value = "'2026-02-05"
$ws.Range("AB19").Value = "'14:56"
$ws.Range("AD19").Value = $False
$ws.Range("AE19").Value = $False
$ws.Range("AG19").Value = $False
$ws.Range("AW19").Value = "Philipp Weiss"
$ws.Range("AX19").Value = "Philipp Weiss"
$ws.Range("M19").ClearContents()
$ws.Range("AC19").ClearContents()

# Row 20
$ws.Range("A20").Value = 131047034
$ws.Range("B20").Value = 78646
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6437
$ws.Range("F20").Value = "Blanksvart spiklav"
$ws.Range("G20").Value = "Calicium denigratum"
$ws.Range("H20").Value = "(Vain.) Tibell"
$ws.Range("P20").Value = "Blomkällan, Dlr"
$ws.Range("Q20").Value = 401597
$ws.Range("R20").Value = 6817852
$ws.Range("S20").Value = 10
$ws.Range("T20").Value = "Dalarna"
$ws.Range("U20").Value = "Älvdalen"
$ws.Range("V20").Value = "Dalarna"
$ws.Range("W20").Value = "Särna"
$ws.Range("Y20").Value = "'2026-02-05"
$ws.Range("Z20").Value = "'15:05"
$ws.Range("AA20").Value = "'2026-02-05"
$ws.Range("AB20").Value = "'15:05"
$ws.Range("AD20").Value = $False
$ws.Range("AE20").Value = $False
$ws.Range("AG20").Value = $False
$ws.Range("AW20").Value = "Philipp Weiss"
$ws.Range("AX20").Value = "Philipp Weiss"
$ws.Range("M20").ClearContents()
$ws.Range("AC20").ClearContents()

# Row 21
$ws.Range("A21").Value = 131046799
$ws.Range("B21").Value = 78255
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 228579
$ws.Range("F21").Value = "Liten svartspik"
$ws.Range("G21").Value = "Chaenothecopsis nana"
$ws.Range("H21").Value = "Tibell"
$ws.Range("P21").Value = "Blomkällan, Dlr"
$ws.Range("Q21").Value = 401649
$ws.Range("R21").Value = 6818014
$ws.Range("S21").Value = 10
$ws.Range("T21").Value = "Dalarna"
$ws.Range("U21").Value = "Älvdalen"
$ws.Range("V21").Value = "Dalarna"
$ws.Range("W21").Value = "Särna"
$ws.Range("Y21").Value = "'2026-02-05"
$ws.Range("Z21").Value = "'14:52"
$ws.Range("AA21").Value = "'2026-02-05"
$ws.Range("AB21").Value = "'14:52"
$ws.Range("AD21").Value = $False
$ws.Range("AE21").Value = $False
$ws.Range("AG21").Value = $False
$ws.Range("AW21").Value = "Philipp Weiss"
$ws.Range("AX21").Value = "Philipp Weiss"

# Row 22
$ws.Range("A22").Value = 131046766
$ws.Range("B22").Value = 57884
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = "Tretåig hackspett"
$ws.Range("G22").Value = "Picoides tridactylus"
$ws.Range("H22").Value = "(Linnaeus, 1758)"
$ws.Range("M22").Value = "äldre spår"
$ws.Range("P22").Value = "Blomkällan, Dlr"
$ws.Range("Q22").Value = 401318
$ws.Range("R22").Value = 6818379
$ws.Range("S22").Value = 10
$ws.Range("T22").Value = "Dalarna"
$ws.Range("U22").Value = "Älvdalen"
$ws.Range("V22").Value = "Dalarna"
$ws.Range("W22").Value = "Särna"
$ws.Range("Y22").Value = "'2026-02-05"
$ws.Range("Z22").Value = "'15:29"
$ws.Range("AA22").Value = "'2026-02-05"
$ws.Range("AB22").Value = "'15:29"
$ws.Range("AC22").Value = "Äldre ringhack (tall)"
$ws.Range("AD22").Value = $False
$ws.Range("AE22").Value = $False
$ws.Range("AG22").Value = $False
$ws.Range("AW22").Value = "Philipp Weiss"
$ws.Range("AX22").Value = "Philipp Weiss"

# Row 25
$ws.Range("A25").Value = 131047014
$ws.Range("B25").Value = 57884
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 100109
$ws.Range("F25").Value = "Tretåig hackspett"
$ws.Range("G25").Value = "Picoides tridactylus"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
$ws.Range("M25").Value = "färska spår"
$ws.Range("P25").Value = "Blomkällan, Dlr"
$ws.Range("Q25").Value = 401378
$ws.Range("R25").Value = 6818082
$ws.Range("S25").Value = 10
$ws.Range("T25").Value = "Dalarna"
$ws.Range("U25").Value = "Älvdalen"
$ws.Range("V25").Value = "Dalarna"
$ws.Range("W25").Value = "Särna"
$ws.Range("Y25").Value = "'2026-02-05"
$ws.Range("Z25").Value = "'15:21"
$ws.Range("AA25").Value = "'2026-02-05"
$ws.Range("AB25").Value = "'15:21"
$ws.Range("AC25").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AD25").Value = $False
$ws.Range("AE25").Value = $True
$ws.Range("AG25").Value = $False
$ws.Range("AW25").Value = "Philipp Weiss"
$ws.Range("AX25").Value = "Philipp Weiss"

# Row 26
$ws.Range("A26").Value = 131046832
$ws.Range("B26").Value = 79243
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("P26").Value = "Blomkällan, Dlr"
$ws.Range("Q26").Value = 401350
$ws.Range("R26").Value = 6818162
$ws.Range("S26").Value = 10
$ws.Range("T26").Value = "Dalarna"
$ws.Range("U26").Value = "Älvdalen"
$ws.Range("V26").Value = "Dalarna"
$ws.Range("W26").Value = "Särna"
$ws.Range("Y26").Value = "'2026-02-05"
$ws.Range("Z26").Value = "'15:24"
$ws.Range("AA26").Value = "'2026-02-05"
$ws.Range("AB26").Value = "'15:24"
$ws.Range("AD26").Value = $False
$ws.Range("AE26").Value = $False
$ws.Range("AG26").Value = $False
$ws.Range("AW26").Value = "Philipp Weiss"
$ws.Range("AX26").Value = "Philipp Weiss"

# Row 27
$ws.Range("A27").Value = 131046826
$ws.Range("B27").Value = 79243
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("P27").Value = "Blomkällan, Dlr"
$ws.Range("Q27").Value = 401647
$ws.Range("R27").Value = 6817965
$ws.Range("S27").Value = 10
$ws.Range("T27").Value = "Dalarna"
$ws.Range("U27").Value = "Älvdalen"
$ws.Range("V27").Value = "Dalarna"
$ws.Range("W27").Value = "Särna"
$ws.Range("Y27").Value = "'2026-02-05"
$ws.Range("Z27").Value = "'14:54"
$ws.Range("AA27").Value = "'2026-02-05"
$ws.Range("AB27").Value = "'14:54"
$ws.Range("AD27").Value = $False
$ws.Range("AE27").Value = $False
$ws.Range("AG27").Value = $False
$ws.Range("AW27").Value = "Philipp Weiss"
$ws.Range("AX27").Value = "Philipp Weiss"
$ws.Range("M27").ClearContents()
$ws.Range("AC27").ClearContents()
